$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Estado de Cuenta" (EC) database rows 16-22: reorder the mora
# periods into chronological order (2305 -> 2311), and refresh the
# "Valor Mora" (F) / "Salario Basico" (G) figures to match the new
# period assignment. Salario Basico drops from 1,160,000 to 950,000 for
# this first batch of the new statement.

$periods = @("2305", "2306", "2307", "2308", "2309", "2310", "2311")
$valorMora = @(46400, 46400, 46400, 46400, 46400, 38000, 35467)
$salarioBasico = 950000

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valorMora[$i]
    $ws.Range("G$row").Value = $salarioBasico
}
